# Refresh per-coin Price (D) / Volume(1h) (E) cells, and the two rows
# whose ranking swapped (EnergySwap <-> Cosmos) - Name (B) and Link (C).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D='66.283.00'; E='  +0.34%  ' }
    @{ Row=3; D='3.574.91'; E='  +2.50%  ' }
    @{ Row=4; E='  +0.00%  ' }
    @{ Row=5; D='609.11' }
    @{ Row=6; D='145.59'; E='  +1.63%  ' }
    @{ Row=7; D='3.574.36'; E='  +2.53%  ' }
    @{ Row=8; E='  +0.21%  ' }
    @{ Row=9; D='0.492'; E='  +3.69%  ' }
    @{ Row=10; E='  +1.33%  ' }
    @{ Row=11; D='7.92'; E='  -3.30%  ' }
    @{ Row=12; E='  +0.59%  ' }
    @{ Row=13; D='4.178.60'; E='  +2.51%  ' }
    @{ Row=14; E='  +2.57%  ' }
    @{ Row=15; D='29.98'; E='  -1.31%  ' }
    @{ Row=16; D='3.566.41'; E='  +2.33%  ' }
    @{ Row=17; D='66.369.61' }
    @{ Row=18; E='  -0.95%  ' }
    @{ Row=19; D='11.53'; E='  +10.99%  ' }
    @{ Row=20; E='  +1.17%  ' }
    @{ Row=21; D='14.90' }
    @{ Row=22; D='430.00'; E='  +2.16%  ' }
    @{ Row=23; D='0.618'; E='  +4.62%  ' }
    @{ Row=24; D='79.21'; E='  +2.22%  ' }
    @{ Row=25; D='3.715.65'; E='  +2.69%  ' }
    @{ Row=26; E='  +0.05%  ' }
    @{ Row=27; E='  +3.65%  ' }
    @{ Row=28; D='2.52'; E='  +2.31%  ' }
    @{ Row=29; D='7.96'; E='  -0.29%  ' }
    @{ Row=30; D='9.10'; E='  -2.50%  ' }
    @{ Row=31; E='  -0.04%  ' }
    @{ Row=32; D='25.68'; E='  +1.92%  ' }
    @{ Row=33; E='  -1.40%  ' }
    @{ Row=34; D='3.568.03'; E='  +2.44%  ' }
    @{ Row=35; D='0.154'; E='  -5.49%  ' }
    @{ Row=37; E='  +1.60%  ' }
    @{ Row=38; D='7.88'; E='  +2.61%  ' }
    @{ Row=39; E='  +0.81%  ' }
    @{ Row=40; D='177.47'; E='  +4.13%  ' }
    @{ Row=41; E='  +0.05%  ' }
    @{ Row=42; E='  -1.60%  ' }
    @{ Row=43; E='  +2.71%  ' }
    @{ Row=44; E='  +0.87%  ' }
    @{ Row=45; D='1.95'; E='  +1.27%  ' }
    @{ Row=46; D='46.18'; E='  +2.44%  ' }
    @{ Row=47; E='  +1.27%  ' }
    @{ Row=48; D='25.69'; E='  -1.48%  ' }
    @{ Row=49; D='2.42'; E='  +2.99%  ' }
    @{ Row=50; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='7.16'; E='  +0.46%  ' }
    @{ Row=51; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='23.54'; E='  +9.31%  ' }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey('B')) { $ws.Cells.Item($r, 2).Value = $u.B }
    if ($u.ContainsKey('C')) { $ws.Cells.Item($r, 3).Value = $u.C }
    if ($u.ContainsKey('D')) {
        $cell = $ws.Cells.Item($r, 4)
        # The feed stores every Price value as literal text. Force Text
        # format first so values that look like plain numbers (e.g.
        # "609.11") are not silently recast as numerics by Excel's
        # automatic type detection on assignment.
        if ($u.D -match '^[0-9]*\.?[0-9]+$') { $cell.NumberFormat = '@' }
        $cell.Value = $u.D
    }
    if ($u.ContainsKey('E')) { $ws.Cells.Item($r, 5).Value = $u.E }
}

